$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new "Comment" column
$ws.Range("D1").Value = "Comment"

# Insert two new rows for UT_5023c (row 26) and UT_5023d (row 27).
# Current layout has row 24 = 5023a, row 25 = 5023b, both referring to
# "Test close loop MUX SQUID". We need to:
#   row24: 5023a | Test close loop MUX SQUID (make pulse)      | NOT PASS | Flux jump due to a bad operating point
#   row25: 5023b | Test close loop MUX SQUID (make pulse)      | PASS
#   row26: 5023c | Test close loop MUX SQUID (make pulseS)     | NOT PASS | Some pixels didn't lock due to a pb in the dmx firmware (feedback integrator incorrectly initialized)
#   row27: 5023d | Test close loop MUX SQUID (make pulseS and auto-relock) | PASS
# then the rest of the rows (previously 26-31) shift down by 2 (to 28-33)

$ws.Rows("26:27").Insert()

$ws.Range("B24").Value = "Test close loop MUX SQUID (make pulse)"
$ws.Range("C24").Value = "NOT PASS"
$ws.Range("D24").Value = "Flux jump due to a bad operating point"

$ws.Range("B25").Value = "Test close loop MUX SQUID (make pulse)"

$ws.Range("A26").Value = "DRE_DMX_UT_5023c"
$ws.Range("B26").Value = "Test close loop MUX SQUID (make pulseS)"
$ws.Range("C26").Value = "NOT PASS"
$ws.Range("D26").Value = "Some pixels didn't lock due to a pb in the dmx firmware (feedback integrator incorrectly initialized)"

$ws.Range("A27").Value = "DRE_DMX_UT_5023d"
$ws.Range("B27").Value = "Test close loop MUX SQUID (make pulseS and auto-relock)"
$ws.Range("C27").Value = "PASS"

# Adjust column widths: column B widened, no longer bestFit; new column D width set
# (target xml widths are 56.42578125 / 34.7109375; the engine quantizes ColumnWidth to
# 1/6-character steps, so these inputs land on the closest achievable values)
$ws.Columns("B").ColumnWidth = 55.67
$ws.Columns("D").ColumnWidth = 33.83

# Set selection / active cell to reflect final state (row after last data row)
$ws.Range("D34").Select()
